$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("PayNowCC")
$ws1.Range("B2").Value = "Sat Nov 15 20:34:58 EST 2025"
$ws1.Range("B3").Value = "Sat Nov 15 20:35:35 EST 2025"
$ws1.Range("B4").Value = "Sat Nov 15 20:36:06 EST 2025"
$ws1.Range("B5").Value = "Sat Nov 15 20:36:35 EST 2025"

$ws2 = $wb.Worksheets.Item("PayNowCCSCF")
$ws2.Range("B2").Value = "Sat Nov 15 20:37:05 EST 2025"
$ws2.Range("B3").Value = "Sat Nov 15 20:37:45 EST 2025"
$ws2.Range("B4").Value = "Sat Nov 15 20:38:25 EST 2025"
$ws2.Range("B5").Value = "Sat Nov 15 20:39:05 EST 2025"

$ws3 = $wb.Worksheets.Item("PayNowCCDCF")
$ws3.Range("B2").Value = "Sat Nov 15 20:39:46 EST 2025"
$ws3.Range("B3").Value = "Sat Nov 15 20:40:25 EST 2025"
$ws3.Range("B4").Value = "Sat Nov 15 20:41:06 EST 2025"
$ws3.Range("B5").Value = "Sat Nov 15 20:41:47 EST 2025"

$ws5 = $wb.Worksheets.Item("OverUnderPay")
$ws5.Range("B2").Value = "Sat Nov 15 20:43:18 EST 2025"
$ws5.Range("B3").Value = "Sat Nov 15 20:46:07 EST 2025"

$ws8 = $wb.Worksheets.Item("NoModifyAmount")
$ws8.Range("B2").Value = "Sat Nov 15 20:46:26 EST 2025"

$ws9 = $wb.Worksheets.Item("NoOverPay")
$ws9.Range("B2").Value = "Sat Nov 15 20:48:15 EST 2025"
